$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the EFA factor correlation matrix (symmetric) with new values
$ws.Range("C2").Value = 0.008
$ws.Range("D2").Value = 0.002
$ws.Range("E2").Value = -0.029

$ws.Range("B3").Value = 0.008
$ws.Range("D3").Value = 0.674
$ws.Range("E3").Value = 0.18

$ws.Range("B4").Value = 0.002
$ws.Range("C4").Value = 0.674
$ws.Range("E4").Value = 0.218

$ws.Range("B5").Value = -0.029
$ws.Range("C5").Value = 0.18
$ws.Range("D5").Value = 0.218
